$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 532; existing rows 532-576 shift down to 533-577.
$ws.Rows(532).Insert()

# Populate the newly inserted row 532 with its full record (mirrors the
# constant columns shared by every data row, plus the row-specific values).
$ws.Cells.Item(532, 1).Value = 3
$ws.Cells.Item(532, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(532, 3).Value = "Coquimbo"
$ws.Cells.Item(532, 4).Value = 45106
$ws.Cells.Item(532, 5).Value = 5
$ws.Cells.Item(532, 6).Value = 100114013
$ws.Cells.Item(532, 7).Value = "Zanahoria"
$ws.Cells.Item(532, 8).Value = "Sin especificar"
$ws.Cells.Item(532, 9).Value = "Primera"
$ws.Cells.Item(532, 10).Value = 230
$ws.Cells.Item(532, 11).Value = 7000
$ws.Cells.Item(532, 12).Value = 7500
$ws.Cells.Item(532, 13).Value = 7239
$ws.Cells.Item(532, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(532, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(532, 16).Value = 362
$ws.Cells.Item(532, 17).Value = 20
$ws.Cells.Item(532, 18).Value = "Hortaliza"
